$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.727.85"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").Value = "'1.862.49"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("D4").Value = "'1.036"
$ws.Range("E4").Value = "  +1.45%  "
$ws.Range("D5").Value = "'323.40"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("D7").Value = "'0.4418"
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("D8").Value = "'0.3793"
$ws.Range("E8").Value = "  +2.13%  "
$ws.Range("D9").Value = "'0.07462"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").Value = "'0.8836"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").Value = "'21.72"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("D12").Value = "'1.868.94"
$ws.Range("E12").Value = "  -9.34%  "
$ws.Range("D13").Value = "'5.542"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "'6.759"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "'0.07220"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "'84.38"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'1.033"
$ws.Range("E19").Value = "  +1.52%  "
$ws.Range("D20").Value = "'15.55"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "'27.729.73"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "'5.305"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").Value = "'11.32"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "'2.090.84"
$ws.Range("E24").Value = "  -7.96%  "
$ws.Range("D25").Value = "'2.016"
$ws.Range("E25").Value = "  +6.48%  "
$ws.Range("D26").Value = "'158.58"
$ws.Range("D27").Value = "'18.83"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").Value = "'1.988"
$ws.Range("E28").Value = "  +3.21%  "
$ws.Range("D29").Value = "'5.324"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").Value = "'118.05"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("D31").Value = "'0.09038"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("D32").Value = "'0.7784"
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("D33").Value = "'1.218"
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("D34").Value = "'3.021"
$ws.Range("E34").Value = "  +6.48%  "
$ws.Range("D35").Value = "'4.573"
$ws.Range("E35").Value = "  +2.59%  "
$ws.Range("E36").Value = "  +1.39%  "
$ws.Range("D37").Value = "'1.151"
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("D38").Value = "'0.01992"
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("D39").Value = "'0.05341"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "'2.870"
$ws.Range("E40").Value = "  +3.09%  "
$ws.Range("D41").Value = "'0.5199"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").Value = "'0.1692"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").Value = "'6.858"
$ws.Range("E43").Value = "  +5.31%  "
$ws.Range("D44").Value = "'8.667"
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("D45").Value = "'110.32"
$ws.Range("E45").Value = "  +2.40%  "
$ws.Range("D46").Value = "'10.63"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("D47").Value = "'0.06630"
$ws.Range("E47").Value = "  +5.70%  "
$ws.Range("D48").Value = "'1.712"
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("D49").Value = "'0.4712"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").Value = "'1.906"
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "'39.73"
$ws.Range("E51").Value = "  +1.58%  "
